$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update rows 2-5 with new order/sum/to_cred/from_merch figures ---
$ws.Range("A2").Value = 2251357
$ws.Range("B2").Value = 20738
$ws.Range("C2").Value = 1.5
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 2.5

$ws.Range("A3").Value = 2250143
$ws.Range("B3").Value = 48178
$ws.Range("C3").Value = 1.5
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 2.5

$ws.Range("A4").Value = 2249706
$ws.Range("B4").Value = 98411
$ws.Range("C4").Value = 1.5
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 2.5

$ws.Range("A5").Value = 2249281
$ws.Range("B5").Value = 17795
$ws.Range("C5").Value = 1.5
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 2.5

# --- Rows 6-11 no longer hold data: wipe the old rows entirely ---
$ws.Range("A6:D11").Style = "Normal"
$ws.Range("A6:D11").ClearContents()
$ws.Range("E6:E11").ClearContents()
$ws.Range("E6:E11").Style = "Normal"

# --- Move the active selection to C7 ---
$ws.Range("C7").Select()
